# Insert a new data row at row 322 (pushes existing rows 322:400 down to 323:401)
# and populate it with the new record described by the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(322).Insert()

$ws.Cells.Item(322, 1).Value  = 4
$ws.Cells.Item(322, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(322, 3).Value  = 'Los Lagos'
$ws.Cells.Item(322, 4).Value  = 44722
$ws.Cells.Item(322, 5).Value  = 10
$ws.Cells.Item(322, 6).Value  = 100114001
$ws.Cells.Item(322, 7).Value  = 'Papa'
$ws.Cells.Item(322, 8).Value  = 'Patagonia'
$ws.Cells.Item(322, 9).Value  = '1a (guarda)'
$ws.Cells.Item(322, 10).Value = 600
$ws.Cells.Item(322, 11).Value = 7000
$ws.Cells.Item(322, 12).Value = 7500
$ws.Cells.Item(322, 13).Value = 7250
$ws.Cells.Item(322, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(322, 15).Value = 'Provincia de Llanquihue'
$ws.Cells.Item(322, 16).Value = 290
$ws.Cells.Item(322, 17).Value = 25
$ws.Cells.Item(322, 18).Value = 'Hortaliza'

# Keep the date style consistent with the rest of column D
$ws.Cells.Item(322, 4).NumberFormat = $ws.Cells.Item(323, 4).NumberFormat
